$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B3 so it is stored as a real number instead of a text value
$ws.Range("B3").Value = 25

# Append the new row of user data (row 4)
$ws.Range("A4").Value = "test"

# B4 must stay a text value ("30"), not be auto-converted to a number,
# so force text entry via a leading apostrophe, then strip the resulting
# quote-prefix formatting so no extra style is applied to the cell.
$ws.Range("B4").Value = "'30"
$ws.Range("B4").ClearFormats()

$ws.Range("C4").Value = "k"
